# Generate Report for Handoff
# The "ae973714-608d-4481-9640-fde0d1bbb16f.md" file has moved from
# "In Translation" to "Ready for handoff" with a fresh handoff timestamp,
# and its priority flips from "ht" (human translation) to "mt" (machine
# translation) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Sheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-02 10:17:53"
$ws.Columns.Item(5).ColumnWidth = 16.25
$ws.Columns.Item(6).ColumnWidth = 16.25

# ---- zh-cn sheet ----
$ws = $wb.Sheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-09-02 10:17:49"
$ws.Columns.Item(3).ColumnWidth = 16.25

# ---- de-de sheet ----
$ws = $wb.Sheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-09-02 10:17:53"
$ws.Columns.Item(3).ColumnWidth = 16.25
